$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 6-10: dates were stored as text strings ("29/03/2020" etc.); convert
# the C column to real date serial values (cell formatting already applies
# the d/mmm/yy display format, so this is a pure value-type change).
$ws.Range("C6").Value = 43919
$ws.Range("C7").Value = 43919
$ws.Range("C8").Value = 43920
$ws.Range("C9").Value = 43920
$ws.Range("C10").Value = 43921

# Row 11 was blank; fill in the new Arithmetic Unit entry for 01/04/2020.
$ws.Range("B11").Value = "`"0624`""
$ws.Range("C11").Value = 43922
$ws.Range("D11").Value = "12:20pm"
$ws.Range("E11").Value = "1:15pm"
$ws.Range("G11").Value = "Working on Arithmetic Unit, trying to fix problems"

# Row 12 was blank; fill in the final entry for 01/04/2020.
$ws.Range("B12").Value = "`"0624`""
$ws.Range("C12").Value = 43922
$ws.Range("D12").Value = "1:30pm"
$ws.Range("E12").Value = "4:17pm"
$ws.Range("G12").Value = "Finished Arithmetic Unit and Logic Unit"

# Update the active cell selection to reflect the last-edited cell.
$ws.Range("E12").Select()
